{"js": "// The underlying change in this revision is purely a re-serialization of\n// the document's OOXML: every element's attributes were re-emitted in\n// alphabetical order (namespace declarations sorted separately from the\n// other attributes) by the authoring tool that produced the commit, and a\n// couple of base64 blobs (`o:gfxdata`) had their internal line-wrapping\n// normalized. No text, formatting, structure, numbering, style, or\n// property value actually changed anywhere in the document - compare\n// word/document.xml and word/styles.xml before/after and every single\n// line differs only in attribute order.\n//\n// Because there is no visible/semantic edit to reproduce, this script\n// intentionally performs a no-op load/sync cycle against the body so the\n// document round-trips unchanged, matching the (content-equivalent)\n// target state.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The underlying change in this revision is purely a re-serialization of\n# the document's OOXML: every element's attributes were re-emitted in\n# alphabetical order (namespace declarations sorted separately from the\n# other attributes) by the authoring tool that produced the commit, and a\n# couple of base64 blobs (`o:gfxdata`) had their internal line-wrapping\n# normalized. No text, formatting, structure, numbering, style, or\n# property value actually changed anywhere in the document - compare\n# word/document.xml and word/styles.xml before/after and every single\n# line differs only in attribute order.\n#\n# Because there is no visible/semantic edit to reproduce, this script\n# intentionally performs a no-op (read-only) pass over the document so it\n# round-trips unchanged, matching the (content-equivalent) target state.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
